$p = $ppt.ActivePresentation

# Slide 7: "OUR SOLUTION " + "AND ITS VALUE  " + "PROPOSITION" (3 runs)
# were merged by the author back into a single run of text while editing
# the title. Re-typing the same text via the Paragraph's TextRange merges
# the split runs into one (keeping the first run's formatting), matching
# how PowerPoint consolidates runs after a user edits text in place.
$slide7 = $p.Slides.Item(7)
$title7 = $slide7.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1)
$title7.Text = "OUR SOLUTION AND ITS VALUE  PROPOSITION TEMP"
$title7.Text = "OUR SOLUTION AND ITS VALUE  PROPOSITION"

# Slide 8: "Dataset " + "description" (2 runs) were merged into one run.
$slide8 = $p.Slides.Item(8)
$title8 = $slide8.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1)
$title8.Text = "Dataset description TEMP"
$title8.Text = "Dataset description"
